$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices, volume %, and two row-position swaps)
# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "0.998") are preserved verbatim instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.650.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.087.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.29"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.080.97"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.94"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.587.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.697.90"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.080.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "504.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.51"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.26"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.24"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "528.32"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -12.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.38"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.82"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.12"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0409"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0790"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.040.51"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.04"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -11.28%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.248"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.95"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.106"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.78"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0490"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -8.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.06%  "
